# Update the FedEx "ShipmentTracking" numbers (column P) for rows 2-16
# with the new tracking numbers issued after the login-email change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTracking = @(
    "320018680212",
    "320018680223",
    "320018680256",
    "320018680278",
    "320018680315",
    "320018680337",
    "320018680360",
    "320018680381",
    "320018680418",
    "320018680430",
    "320018680473",
    "320018680495",
    "320018680521",
    "320018680543",
    "320018680576"
)

# Tracking numbers are long digit strings that Excel would otherwise parse
# as numeric values. Force the target range to Text format first so the
# values are stored as text (matching the original file's shared-string
# cells), then clear the temporary formatting so the cells keep their
# original (unstyled) look.
$targetRange = $ws.Range("P2:P16")
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $newTracking.Count; $i++) {
    $row = 2 + $i
    $ws.Range("P" + $row).Value = $newTracking[$i]
}

$targetRange.ClearFormats()
